$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added the 2nd reset pin -> need a 2nd unit of the 100nF capacitor (row 9)
# Quantity 1 -> 2, and the line cost doubles 0.35 -> 0.7 accordingly.
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0.7

# Leave the cursor where the author left it when they saved.
$ws.Range("H10").Select()
